$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.747.02"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "3.467.95"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "581.56"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "167.91"
$ws.Range("E6").Value = "  +4.43%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.467.85"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "0.564"
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("D10").Value = "7.31"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").Value = "0.123"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "0.431"
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("D13").Value = "4.051.53"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "27.54"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "0.0000176"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "65.685.46"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "3.444.46"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").Value = "6.24"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "13.78"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "384.31"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "7.96"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "71.76"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("D25").Value = "0.521"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("D26").Value = "'0.0000120"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").Value = "9.85"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "0.182"
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "6.24"
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D33").Value = "23.32"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "7.32"
$ws.Range("E34").Value = "  +3.58%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "1.53"
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("D37").Value = "160.04"
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").Value = "0.891"
$ws.Range("E38").Value = "  +8.33%  "
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "0.0737"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "6.65"
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "26.26"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.810.73"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "26.76"
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "43.17"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").Value = "4.47"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0311"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "2.49"
$ws.Range("E48").Value = "  +2.63%  "
$ws.Range("D49").Value = "338.41"
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("D50").Value = "1.07"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "32.44"
$ws.Range("E51").Value = "  +4.82%  "
